$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E17").Value = "1804"
$ws.Range("F17").Value = 31249

$ws.Range("E18").Value = "1803"
$ws.Range("F18").Value = 31249

$ws.Range("E19").Value = "1802"
$ws.Range("F19").Value = 31249

$ws.Range("E20").Value = "1801"
$ws.Range("F20").Value = 16666
